$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 67

$ws.Cells.Item($row, 1).Value = "Kindergarden"
$ws.Cells.Item($row, 2).Value = "Kindergarden Delft Hugo de Grootstraat"
$ws.Cells.Item($row, 3).Value = "BSO"

# Column D holds a date-shaped string ("2024-02-13") that must stay literal
# text (matching the other short-form dates already in the sheet), not get
# auto-converted into a date serial number. Writing it through a formula
# that evaluates to a text string, then collapsing it to a static value via
# copy / paste-special-values, avoids Excel's text->date autodetection
# (which a direct .Value assignment would trigger) without touching
# NumberFormat/Style (which would otherwise leave a stray new style behind).
$ws.Cells.Item($row, 4).Formula = "=""2024-02-13"""
$ws.Cells.Item($row, 4).Copy()
$ws.Cells.Item($row, 4).PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 0
